$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "./images3/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_Ref_BG-grey_stim-white.png",
    "./images3/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_Ref_BG-grey_stim-white.png",
    "./images3/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_Ref_BG-grey_stim-white.png",
    "./images3/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_Ref_BG-grey_stim-white.png",
    "./images3/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images3/Sphere_CW-3.75_BG-grey_stim-white.png"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
